$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("B2").Value = 0.171875
$ws.Range("C2").Value = 0.6041666666666666
$ws.Range("J2").Value = 0.005208333333333333
$ws.Range("P2").Value = 0.140625
$ws.Range("S2").Value = 0.078125
$ws.Range("B3").Value = 0.02419354838709677
$ws.Range("C3").Value = 0.0564516129032258
$ws.Range("J3").Value = 0.008064516129032258
$ws.Range("P3").Value = 0.7016129032258065
$ws.Range("S3").Value = 0.2096774193548387
$ws.Range("J4").Value = 0.02272727272727273
$ws.Range("P4").Value = 0.7954545454545454
$ws.Range("S4").Value = 0.1818181818181818
$ws.Range("B6").Value = 0.0374331550802139
$ws.Range("D6").Value = 0.0160427807486631
$ws.Range("E6").Value = 0.0106951871657754
$ws.Range("F6").Value = 0.09090909090909091
$ws.Range("J6").Value = 0.2192513368983957
$ws.Range("O6").Value = 0.0160427807486631
$ws.Range("Q6").Value = 0.1711229946524064
$ws.Range("R6").Value = 0.09090909090909091
$ws.Range("S6").Value = 0.3475935828877005
$ws.Range("B7").Value = 0.07462686567164178
$ws.Range("D7").Value = 0.004975124378109453
$ws.Range("F7").Value = 0.04477611940298507
$ws.Range("J7").Value = 0.1293532338308458
$ws.Range("O7").Value = 0.009950248756218905
$ws.Range("Q7").Value = 0.1791044776119403
$ws.Range("R7").Value = 0.06965174129353234
$ws.Range("S7").Value = 0.4875621890547264
$ws.Range("B8").Value = 0.06199460916442048
$ws.Range("D8").Value = 0.01886792452830189
$ws.Range("F8").Value = 0.0431266846361186
$ws.Range("J8").Value = 0.1266846361185984
$ws.Range("O8").Value = 0.0215633423180593
$ws.Range("Q8").Value = 0.1644204851752022
$ws.Range("R8").Value = 0.1078167115902965
$ws.Range("S8").Value = 0.4555256064690027
$ws.Range("B9").Value = 0.06511627906976744
$ws.Range("D9").Value = 0.01395348837209302
$ws.Range("F9").Value = 0.06046511627906977
$ws.Range("J9").Value = 0.1441860465116279
$ws.Range("O9").Value = 0.02790697674418605
$ws.Range("Q9").Value = 0.213953488372093
$ws.Range("R9").Value = 0.09302325581395349
$ws.Range("S9").Value = 0.3813953488372093
$ws.Range("B10").Value = 0.0815485996705107
$ws.Range("D10").Value = 0.02800658978583196
$ws.Range("F10").Value = 0.0700164744645799
$ws.Range("J10").Value = 0.1383855024711697
$ws.Range("O10").Value = 0.009884678747940691
$ws.Range("Q10").Value = 0.2075782537067545
$ws.Range("R10").Value = 0.09637561779242175
$ws.Range("S10").Value = 0.3682042833607908
$ws.Range("G11").Value = 0.1742160278745645
$ws.Range("J11").Value = 0.04181184668989547
$ws.Range("K11").Value = 0.2090592334494774
$ws.Range("L11").Value = 0.554006968641115
$ws.Range("S11").Value = 0.02090592334494774
$ws.Range("G12").Value = 0.7590361445783133
$ws.Range("J12").Value = 0.1746987951807229
$ws.Range("K12").Value = 0.006024096385542169
$ws.Range("L12").Value = 0.01807228915662651
$ws.Range("S12").Value = 0.04216867469879518
$ws.Range("G13").Value = 0.6511627906976745
$ws.Range("J13").Value = 0.3023255813953488
$ws.Range("S13").Value = 0.04651162790697674
$ws.Range("G14").Value = 0.5
$ws.Range("J14").Value = 0.5
$ws.Range("F15").Value = 0.01212121212121212
$ws.Range("H15").Value = 0.1575757575757576
$ws.Range("I15").Value = 0.08484848484848485
$ws.Range("J15").Value = 0.3757575757575757
$ws.Range("K15").Value = 0.04848484848484848
$ws.Range("M15").Value = 0.006060606060606061
$ws.Range("O15").Value = 0.04242424242424243
$ws.Range("S15").Value = 0.2727272727272727
$ws.Range("F16").Value = 0.00684931506849315
$ws.Range("H16").Value = 0.1712328767123288
$ws.Range("I16").Value = 0.1301369863013699
$ws.Range("J16").Value = 0.4178082191780822
$ws.Range("K16").Value = 0.1164383561643836
$ws.Range("M16").Value = 0.03424657534246575
$ws.Range("O16").Value = 0.0273972602739726
$ws.Range("S16").Value = 0.0958904109589041
$ws.Range("F17").Value = 0.007125890736342043
$ws.Range("H17").Value = 0.159144893111639
$ws.Range("I17").Value = 0.1211401425178147
$ws.Range("J17").Value = 0.4180522565320665
$ws.Range("K17").Value = 0.1187648456057007
$ws.Range("M17").Value = 0.01187648456057007
$ws.Range("O17").Value = 0.05463182897862233
$ws.Range("S17").Value = 0.1092636579572447
$ws.Range("F18").Value = 0.01435406698564593
$ws.Range("H18").Value = 0.1913875598086124
$ws.Range("I18").Value = 0.1004784688995215
$ws.Range("J18").Value = 0.430622009569378
$ws.Range("K18").Value = 0.1100478468899522
$ws.Range("M18").Value = 0.009569377990430622
$ws.Range("N18").Value = 0.004784688995215311
$ws.Range("O18").Value = 0.04784688995215311
$ws.Range("S18").Value = 0.09090909090909091
$ws.Range("F19").Value = 0.01757469244288225
$ws.Range("H19").Value = 0.1880492091388401
$ws.Range("I19").Value = 0.09578207381370826
$ws.Range("J19").Value = 0.4103690685413005
$ws.Range("K19").Value = 0.1107205623901582
$ws.Range("M19").Value = 0.02636203866432337
$ws.Range("N19").Value = 0.001757469244288225
$ws.Range("O19").Value = 0.05448154657293498
$ws.Range("S19").Value = 0.09490333919156414
